# backlog.xlsx - add "Tekst" header / "Totaal" sum row, turn the range into
# a formatted table ("Tabel4"), and refresh the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the stray blank row 19, pulling the SUM row up to row 19 ---
$ws.Rows(19).Delete()

# --- 2. New header cell for column D ---
$ws.Range("D2").Value = "Tekst"

# --- 3. New "Totaal" label + keep the SUM formula on row 19 ---
$ws.Range("D19").Value = "Totaal"
$ws.Range("E19").Formula = "=SUM(E3:E18)"

# --- 4. Cell-level formatting (font + alignment) ---
# Header row + the numeric "value" column keep their centered look, just on
# the (new) font used throughout the rest of the sheet.
$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("E2").Font.Name = "Calibri"
$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("E3:E18").Font.Name = "Calibri"

# "Volgorde" numbers + the new "Totaal" label -> right aligned
$ws.Range("C3:C18").Font.Name = "Calibri"
$ws.Range("C3:C18").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight
$ws.Range("D19").Font.Name = "Calibri"
$ws.Range("D19").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

# Free text column + "Opmerkingen" cell -> plain font, default alignment
$ws.Range("D2:D18").Font.Name = "Calibri"
$ws.Range("F9").Font.Name = "Calibri"

# SUM total -> left aligned
$ws.Range("E19").Font.Name = "Calibri"
$ws.Range("E19").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft

# --- 5. Turn C2:F19 into a real table ---
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("C2:F19"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Tabel4"
$tbl.TableStyle = "TableStyleMedium7"

# --- 6. Selection, as last seen in the authored workbook ---
$ws.Range("H10").Select()
